# Auto-generated PowerShell Excel COM-interop script
# Applies the cryptos.xlsx update described by the commit diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the Price/Volume columns to Text first so that numeric-looking
# strings (e.g. "1.003") are stored as text (matching the original
# inline-string cells) instead of being auto-converted to numbers.
$dataRange = $ws.Range("D2:E51")
$dataRange.NumberFormat = "@"

# Row 2
$ws.Range("D2").Value = "27.890.95"
$ws.Range("E2").Value = "  -1.99%  "

# Row 3
$ws.Range("D3").Value = "1.864.57"
$ws.Range("E3").Value = "  -2.65%  "

# Row 4
$ws.Range("D4").Value = "1.003"
$ws.Range("E4").Value = "  -0.02%  "

# Row 5
$ws.Range("D5").Value = "312.65"
$ws.Range("E5").Value = "  -0.98%  "

# Row 6
$ws.Range("E6").Value = "  -0.02%  "

# Row 7
$ws.Range("D7").Value = "0.4991"
$ws.Range("E7").Value = "  -2.33%  "

# Row 8
$ws.Range("D8").Value = "0.3818"
$ws.Range("E8").Value = "  -3.69%  "

# Row 9
$ws.Range("D9").Value = "0.08862"
$ws.Range("E9").Value = "  -8.75%  "

# Row 10
$ws.Range("D10").Value = "1.116"
$ws.Range("E10").Value = "  -2.38%  "

# Row 11
$ws.Range("D11").Value = "41.48"
$ws.Range("E11").Value = "  -1.52%  "

# Row 12
$ws.Range("D12").Value = "6.361"
$ws.Range("E12").Value = "  -1.63%  "

# Row 13
$ws.Range("D13").Value = "20.58"
$ws.Range("E13").Value = "  -2.14%  "

# Row 14
$ws.Range("D14").Value = "1.865.95"
$ws.Range("E14").Value = "  -2.94%  "

# Row 15
$ws.Range("D15").Value = "7.217"
$ws.Range("E15").Value = "  -2.41%  "

# Row 16
$ws.Range("E16").Value = "  +0.08%  "

# Row 17
$ws.Range("D17").Value = "0.00001094"
$ws.Range("E17").Value = "  -3.17%  "

# Row 18
$ws.Range("D18").Value = "90.93"
$ws.Range("E18").Value = "  -3.08%  "

# Row 19
$ws.Range("D19").Value = "0.06670"
$ws.Range("E19").Value = "  -0.11%  "

# Row 20
$ws.Range("D20").Value = "17.93"
$ws.Range("E20").Value = "  -0.86%  "

# Row 21
$ws.Range("D21").Value = "1.001"
$ws.Range("E21").Value = "  -0.05%  "

# Row 22
$ws.Range("D22").Value = "6.092"
$ws.Range("E22").Value = "  -2.73%  "

# Row 23
$ws.Range("D23").Value = "27.932.36"
$ws.Range("E23").Value = "  -2.04%  "

# Row 24
$ws.Range("D24").Value = "11.44"
$ws.Range("E24").Value = "  -0.08%  "

# Row 25
$ws.Range("D25").Value = "2.286"
$ws.Range("E25").Value = "  -1.42%  "

# Row 26
$ws.Range("D26").Value = "2.075.27"
$ws.Range("E26").Value = "  -2.92%  "

# Row 27
$ws.Range("D27").Value = "2.489"
$ws.Range("E27").Value = "  -6.75%  "

# Row 28
$ws.Range("D28").Value = "157.61"
$ws.Range("E28").Value = "  -0.58%  "

# Row 29
$ws.Range("E29").Value = "  -2.58%  "

# Row 30
$ws.Range("D30").Value = "126.31"
$ws.Range("E30").Value = "  -1.61%  "

# Row 31
$ws.Range("D31").Value = "0.1057"
$ws.Range("E31").Value = "  -1.18%  "

# Row 32
$ws.Range("D32").Value = "1.051"
$ws.Range("E32").Value = "  -4.55%  "

# Row 33
$ws.Range("D33").Value = "5.574"
$ws.Range("E33").Value = "  -2.07%  "

# Row 34
$ws.Range("D34").Value = "3.600"
$ws.Range("E34").Value = "  -1.11%  "

# Row 35
$ws.Range("D35").Value = "9.327"
$ws.Range("E35").Value = "  -5.12%  "

# Row 36
$ws.Range("D36").Value = "0.06508"
$ws.Range("E36").Value = "  -3.01%  "

# Row 37
$ws.Range("D37").Value = "0.02389"
$ws.Range("E37").Value = "  -2.07%  "

# Row 38
$ws.Range("D38").Value = "0.2177"
$ws.Range("E38").Value = "  -1.95%  "

# Row 39
$ws.Range("D39").Value = "1.273"
$ws.Range("E39").Value = "  +5.04%  "

# Row 40
$ws.Range("D40").Value = "1.194"
$ws.Range("E40").Value = "  -4.70%  "

# Row 41
$ws.Range("B41").Value = "TheSandbox"
$ws.Range("C41").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D41").Value = "0.6355"
$ws.Range("E41").Value = "  -1.09%  "

# Row 42
$ws.Range("B42").Value = "Aptos"
$ws.Range("C42").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D42").Value = "11.40"
$ws.Range("E42").Value = "  -2.00%  "

# Row 43
$ws.Range("D43").Value = "4.901"
$ws.Range("E43").Value = "  -3.06%  "

# Row 44
$ws.Range("D44").Value = "1.002"
$ws.Range("E44").Value = "  -0.03%  "

# Row 45
$ws.Range("D45").Value = "13.15"
$ws.Range("E45").Value = "  -2.83%  "

# Row 46
$ws.Range("D46").Value = "0.5987"
$ws.Range("E46").Value = "  -1.22%  "

# Row 47
$ws.Range("D47").Value = "1.284"
$ws.Range("E47").Value = "  +0.06%  "

# Row 48
$ws.Range("D48").Value = "3.674"
$ws.Range("E48").Value = "  -2.94%  "

# Row 49
$ws.Range("D49").Value = "1.221"
$ws.Range("E49").Value = "  +2.04%  "

# Row 50
$ws.Range("D50").Value = "1.981"
$ws.Range("E50").Value = "  -3.47%  "

# Row 51
$ws.Range("D51").Value = "120.55"
$ws.Range("E51").Value = "  -3.53%  "

# Restore the default (General) formatting on the Price/Volume columns so
# no stray number-format/style changes are introduced by the edit above.
$dataRange.ClearFormats()
